# The "sold" lookups were timing out before the page finished loading, so
# the stale/sold listings below are pruned and the freshly re-checked
# listings are appended at the bottom of the URL list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows (1-based, matching the sheet's current A2:A29 URL block) that are no
# longer valid and must be removed. Clear every existing hyperlink first so
# the Hyperlinks collection + its relationship ids don't end up stale once
# rows start shifting around.
$ws.Range("A1").Hyperlinks.Delete()

$rowsToDelete = @(28, 19, 16, 13, 8, 3)
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}

# New listings discovered on the slower re-check; appended after the
# surviving rows.
$newUrls = @(
    "https://www.autotrader.co.uk/car-details/202503130120222?fromSavedAds=true&advertising-location=at_cars&sort=relevance&postcode=CB58TJ",
    "https://www.autotrader.co.uk/car-details/202502199274313?fromSavedAds=true&advertising-location=at_cars&sort=relevance&postcode=CB58TJ",
    "https://www.autotrader.co.uk/car-details/202503150168260?fromSavedAds=true&advertising-location=at_cars&sort=relevance&postcode=CB58TJ",
    "https://www.autotrader.co.uk/car-details/202503200354314?fromSavedAds=true&advertising-location=at_cars&sort=relevance&postcode=CB58TJ",
    "https://www.autotrader.co.uk/car-details/202503170247683?fromSavedAds=true&advertising-location=at_cars&sort=relevance&postcode=CB58TJ",
    "https://www.autotrader.co.uk/car-details/202501178145053?fromSavedAds=true&advertising-location=at_cars&sort=relevance&postcode=CB58TJ",
    "https://www.autotrader.co.uk/car-details/202502259471660?fromSavedAds=true&advertising-location=at_cars&sort=relevance&postcode=CB58TJ"
)

$lastRow = $ws.UsedRange.Rows.Count
foreach ($u in $newUrls) {
    $lastRow = $lastRow + 1
    $ws.Cells.Item($lastRow, 1).Value = $u
    $ws.Cells.Item($lastRow, 1).Style = "Hyperlink"
}

# Rebuild the hyperlinks for every URL row (A2 downward) now that the final
# row order is settled.
$finalLastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $finalLastRow; $r++) {
    $target = $ws.Cells.Item($r, 1).Text
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 1), $target, "", "", $target)
    $ws.Cells.Item($r, 1).Style = "Hyperlink"
}
